$d = $word.ActiveDocument

# --- Split the Background run for Alex into three runs ---
$d.Content.Find.Execute(
    ": Alex is an avid gamer and a computer science major",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": Alex is a gamer and a computer science major", 2)

Write-Output "done"
